$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.636.07"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").Value = "'1.662.83"
$ws.Range("E3").Value = "  +3.30%  "

$ws.Range("D4").Value = "'0.9980"
$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").Value = "'0.9988"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").Value = "'302.50"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "'0.3840"
$ws.Range("E7").Value = "  +1.58%  "

$ws.Range("D8").Value = "'0.3608"
$ws.Range("E8").Value = "  +2.38%  "

$ws.Range("D9").Value = "'51.16"
$ws.Range("E9").Value = "  -1.15%  "

$ws.Range("D10").Value = "'1.246"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").Value = "'0.08201"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "'22.50"
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").Value = "'6.525"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "'7.527"
$ws.Range("E15").Value = "  +4.43%  "

$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("D17").Value = "'1.654.39"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").Value = "'97.67"
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("D19").Value = "'0.06995"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("D20").Value = "'6.857"
$ws.Range("E20").Value = "  +5.34%  "

$ws.Range("D21").Value = "'17.78"
$ws.Range("E21").Value = "  +3.84%  "

$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").Value = "'12.74"
$ws.Range("E23").Value = "  +3.62%  "

$ws.Range("D24").Value = "'23.593.65"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("D25").Value = "'2.517"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").Value = "'3.012"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").Value = "'21.25"
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("D28").Value = "'153.28"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").Value = "'5.235"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "'134.27"
$ws.Range("E30").Value = "  +1.58%  "

$ws.Range("D31").Value = "'7.241"
$ws.Range("E31").Value = "  +11.78%  "

$ws.Range("D32").Value = "'1.838.86"
$ws.Range("E32").Value = "  +2.75%  "

$ws.Range("D33").Value = "'2.244"
$ws.Range("E33").Value = "  +7.13%  "

$ws.Range("D34").Value = "'11.96"
$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("D35").Value = "'1.058"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("E36").Value = "  +4.01%  "

$ws.Range("D37").Value = "'6.161"
$ws.Range("E37").Value = "  +6.13%  "

$ws.Range("D38").Value = "'0.2505"
$ws.Range("E38").Value = "  +2.54%  "

$ws.Range("D39").Value = "'0.08804"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").Value = "'0.07048"
$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("D41").Value = "'13.34"
$ws.Range("E41").Value = "  +11.90%  "

$ws.Range("D42").Value = "'0.7041"
$ws.Range("E42").Value = "  +2.93%  "

$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").Value = "'16.14"
$ws.Range("E44").Value = "  +6.33%  "

$ws.Range("D45").Value = "'0.6572"
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.317"
$ws.Range("E46").Value = "  +3.61%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").Value = "'3.960"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "'0.07959"
$ws.Range("E49").Value = "  +1.42%  "

$ws.Range("D50").Value = "'128.17"
$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("D51").Value = "'1.200"
$ws.Range("E51").Value = "  +3.46%  "
